$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Activité des personnes âgées de 15 ans et plus"
$ws.Range("B16").Value = "Taux d’analphabétisme"

[void]$ws.Range("B16").Select()
